# The source diff for this commit ("rebased barcode to main") touches only
# the root-element namespace-declaration ordering in word/document.xml,
# word/footer.xml, word/header.xml and word/styles.xml (e.g. xmlns:r / xmlns:w
# / xmlns:w15 / ... being emitted in a different sequence). Every attribute
# that appears before the change still appears after it (same prefixes, same
# URIs, same mc:Ignorable token list) and every other byte of each part is
# unchanged - this is a cosmetic re-serialization artifact produced by the
# tool that regenerated the fixture (docx4j) during the rebase, not a
# document edit. There is no corresponding Word object-model mutation
# (paragraph/run/style/property change) to make: the document's content,
# formatting, headers and footers are identical before and after.
#
# So this script intentionally performs no content changes - it leaves
# $word.ActiveDocument exactly as loaded.
$d = $word.ActiveDocument
